$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the whole used range one column to the left (B1:F3 -> A1:E3),
# which carries each cell's original formatting (header style, shared
# string refs, values) along with it.
$ws.Range("B1:F3").Copy($ws.Range("A1"))

# The old rightmost column (F) is now a duplicate leftover; remove it
# completely (value + formatting) so the used range shrinks to A1:E3.
$ws.Range("F1:F3").Clear()
